$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# --- Step 1: simple text replacements (cells whose position/row does not change) ---
$d.Content.Find.Execute("13+54=", $true, $false, $false, $false, $false, $true, 1, $false, "18-10=", 2) | Out-Null
$d.Content.Find.Execute("77-42=", $true, $false, $false, $false, $false, $true, 1, $false, "60-46=", 2) | Out-Null
$d.Content.Find.Execute("94-48=", $true, $false, $false, $false, $false, $true, 1, $false, "47+49=", 2) | Out-Null
$d.Content.Find.Execute("14+9=", $true, $false, $false, $false, $false, $true, 1, $false, "61+11=", 2) | Out-Null
$d.Content.Find.Execute("97-70=", $true, $false, $false, $false, $false, $true, 1, $false, "47-3=", 2) | Out-Null
$d.Content.Find.Execute("0+83=", $true, $false, $false, $false, $false, $true, 1, $false, "36+42=", 2) | Out-Null
$d.Content.Find.Execute("16+44=", $true, $false, $false, $false, $false, $true, 1, $false, "39+47=", 2) | Out-Null
$d.Content.Find.Execute("9+36=", $true, $false, $false, $false, $false, $true, 1, $false, "81-12=", 2) | Out-Null
$d.Content.Find.Execute("12+73=", $true, $false, $false, $false, $false, $true, 1, $false, "75-5=", 2) | Out-Null
$d.Content.Find.Execute("7+89=", $true, $false, $false, $false, $false, $true, 1, $false, "87-66=", 2) | Out-Null
$d.Content.Find.Execute("46+49=", $true, $false, $false, $false, $false, $true, 1, $false, "13+31=", 2) | Out-Null
$d.Content.Find.Execute("45-33=", $true, $false, $false, $false, $false, $true, 1, $false, "92-8=", 2) | Out-Null
$d.Content.Find.Execute("22+8=", $true, $false, $false, $false, $false, $true, 1, $false, "28+6=", 2) | Out-Null
$d.Content.Find.Execute("84-22=", $true, $false, $false, $false, $false, $true, 1, $false, "29+52=", 2) | Out-Null
$d.Content.Find.Execute("8+32=", $true, $false, $false, $false, $false, $true, 1, $false, "2+2=", 2) | Out-Null
$d.Content.Find.Execute("71-49=", $true, $false, $false, $false, $false, $true, 1, $false, "17-8=", 2) | Out-Null
$d.Content.Find.Execute("74-36=", $true, $false, $false, $false, $false, $true, 1, $false, "99-97=", 2) | Out-Null
$d.Content.Find.Execute("26+39=", $true, $false, $false, $false, $false, $true, 1, $false, "68-63=", 2) | Out-Null
$d.Content.Find.Execute("77+16=", $true, $false, $false, $false, $false, $true, 1, $false, "79-64=", 2) | Out-Null
$d.Content.Find.Execute("13-10=", $true, $false, $false, $false, $false, $true, 1, $false, "79-38=", 2) | Out-Null
$d.Content.Find.Execute("60-50=", $true, $false, $false, $false, $false, $true, 1, $false, "72+12=", 2) | Out-Null
$d.Content.Find.Execute("35+27=", $true, $false, $false, $false, $false, $true, 1, $false, "16-2=", 2) | Out-Null
$d.Content.Find.Execute("29+15=", $true, $false, $false, $false, $false, $true, 1, $false, "60-47=", 2) | Out-Null
$d.Content.Find.Execute("28-21=", $true, $false, $false, $false, $false, $true, 1, $false, "85-31=", 2) | Out-Null
$d.Content.Find.Execute("32+65=", $true, $false, $false, $false, $false, $true, 1, $false, "16+74=", 2) | Out-Null
$d.Content.Find.Execute("25+45=", $true, $false, $false, $false, $false, $true, 1, $false, "45+1=", 2) | Out-Null
$d.Content.Find.Execute("76+16=", $true, $false, $false, $false, $false, $true, 1, $false, "61-16=", 2) | Out-Null
$d.Content.Find.Execute("74-3=", $true, $false, $false, $false, $false, $true, 1, $false, "65+15=", 2) | Out-Null
$d.Content.Find.Execute("87-77=", $true, $false, $false, $false, $false, $true, 1, $false, "73+7=", 2) | Out-Null
$d.Content.Find.Execute("59+25=", $true, $false, $false, $false, $false, $true, 1, $false, "65+17=", 2) | Out-Null
$d.Content.Find.Execute("16+60=", $true, $false, $false, $false, $false, $true, 1, $false, "24+61=", 2) | Out-Null
$d.Content.Find.Execute("20+26=", $true, $false, $false, $false, $false, $true, 1, $false, "14+81=", 2) | Out-Null
$d.Content.Find.Execute("65+16=", $true, $false, $false, $false, $false, $true, 1, $false, "44-12=", 2) | Out-Null
$d.Content.Find.Execute("23+51=", $true, $false, $false, $false, $false, $true, 1, $false, "94+1=", 2) | Out-Null
$d.Content.Find.Execute("64-26=", $true, $false, $false, $false, $false, $true, 1, $false, "92-76=", 2) | Out-Null
$d.Content.Find.Execute("26+66=", $true, $false, $false, $false, $false, $true, 1, $false, "71-35=", 2) | Out-Null
$d.Content.Find.Execute("33+2=", $true, $false, $false, $false, $false, $true, 1, $false, "35-17=", 2) | Out-Null
$d.Content.Find.Execute("82-17=", $true, $false, $false, $false, $false, $true, 1, $false, "79-42=", 2) | Out-Null
$d.Content.Find.Execute("75+11=", $true, $false, $false, $false, $false, $true, 1, $false, "7+66=", 2) | Out-Null
$d.Content.Find.Execute("56-48=", $true, $false, $false, $false, $false, $true, 1, $false, "13+36=", 2) | Out-Null
$d.Content.Find.Execute("74+7=", $true, $false, $false, $false, $false, $true, 1, $false, "91-28=", 2) | Out-Null
$d.Content.Find.Execute("11+36=", $true, $false, $false, $false, $false, $true, 1, $false, "82-54=", 2) | Out-Null
$d.Content.Find.Execute("62-16=", $true, $false, $false, $false, $false, $true, 1, $false, "25+30=", 2) | Out-Null
$d.Content.Find.Execute("66-14=", $true, $false, $false, $false, $false, $true, 1, $false, "28+23=", 2) | Out-Null
$d.Content.Find.Execute("77-23=", $true, $false, $false, $false, $false, $true, 1, $false, "58+22=", 2) | Out-Null
$d.Content.Find.Execute("10+15=", $true, $false, $false, $false, $false, $true, 1, $false, "29+50=", 2) | Out-Null
$d.Content.Find.Execute("63+32=", $true, $false, $false, $false, $false, $true, 1, $false, "36+49=", 2) | Out-Null
$d.Content.Find.Execute("56-18=", $true, $false, $false, $false, $false, $true, 1, $false, "60+27=", 2) | Out-Null
$d.Content.Find.Execute("17+41=", $true, $false, $false, $false, $false, $true, 1, $false, "19+47=", 2) | Out-Null
$d.Content.Find.Execute("27+54=", $true, $false, $false, $false, $false, $true, 1, $false, "95-41=", 2) | Out-Null
$d.Content.Find.Execute("25+65=", $true, $false, $false, $false, $false, $true, 1, $false, "89-17=", 2) | Out-Null
$d.Content.Find.Execute("37+15=", $true, $false, $false, $false, $false, $true, 1, $false, "39-10=", 2) | Out-Null
$d.Content.Find.Execute("90-43=", $true, $false, $false, $false, $false, $true, 1, $false, "88-16=", 2) | Out-Null
$d.Content.Find.Execute("50-35=", $true, $false, $false, $false, $false, $true, 1, $false, "54+35=", 2) | Out-Null
$d.Content.Find.Execute("61-32=", $true, $false, $false, $false, $false, $true, 1, $false, "67-56=", 2) | Out-Null
$d.Content.Find.Execute("16-14=", $true, $false, $false, $false, $false, $true, 1, $false, "55+35=", 2) | Out-Null
$d.Content.Find.Execute("36+50=", $true, $false, $false, $false, $false, $true, 1, $false, "26-22=", 2) | Out-Null
$d.Content.Find.Execute("44-5=", $true, $false, $false, $false, $false, $true, 1, $false, "77-26=", 2) | Out-Null
$d.Content.Find.Execute("67-7=", $true, $false, $false, $false, $false, $true, 1, $false, "66+16=", 2) | Out-Null
$d.Content.Find.Execute("50+17=", $true, $false, $false, $false, $false, $true, 1, $false, "95-76=", 2) | Out-Null
$d.Content.Find.Execute("22+34=", $true, $false, $false, $false, $false, $true, 1, $false, "52-14=", 2) | Out-Null
$d.Content.Find.Execute("67-36=", $true, $false, $false, $false, $false, $true, 1, $false, "30+23=", 2) | Out-Null
$d.Content.Find.Execute("34-10=", $true, $false, $false, $false, $false, $true, 1, $false, "46+53=", 2) | Out-Null
$d.Content.Find.Execute("61-5=", $true, $false, $false, $false, $false, $true, 1, $false, "25+51=", 2) | Out-Null
$d.Content.Find.Execute("81+17=", $true, $false, $false, $false, $false, $true, 1, $false, "80-0=", 2) | Out-Null
$d.Content.Find.Execute("29-11=", $true, $false, $false, $false, $false, $true, 1, $false, "17+35=", 2) | Out-Null
$d.Content.Find.Execute("22+6=", $true, $false, $false, $false, $false, $true, 1, $false, "95-10=", 2) | Out-Null
$d.Content.Find.Execute("26+56=", $true, $false, $false, $false, $false, $true, 1, $false, "2+68=", 2) | Out-Null
$d.Content.Find.Execute("40+4=", $true, $false, $false, $false, $false, $true, 1, $false, "31+6=", 2) | Out-Null
$d.Content.Find.Execute("85-83=", $true, $false, $false, $false, $false, $true, 1, $false, "27+62=", 2) | Out-Null
$d.Content.Find.Execute("76+19=", $true, $false, $false, $false, $false, $true, 1, $false, "3+78=", 2) | Out-Null
$d.Content.Find.Execute("77-16=", $true, $false, $false, $false, $false, $true, 1, $false, "65+33=", 2) | Out-Null
$d.Content.Find.Execute("29+57=", $true, $false, $false, $false, $false, $true, 1, $false, "76-42=", 2) | Out-Null
$d.Content.Find.Execute("72+10=", $true, $false, $false, $false, $false, $true, 1, $false, "49-16=", 2) | Out-Null
$d.Content.Find.Execute("80+16=", $true, $false, $false, $false, $false, $true, 1, $false, "68-53=", 2) | Out-Null
$d.Content.Find.Execute("72-7=", $true, $false, $false, $false, $false, $true, 1, $false, "72-48=", 2) | Out-Null
$d.Content.Find.Execute("70-38=", $true, $false, $false, $false, $false, $true, 1, $false, "17+46=", 2) | Out-Null
$d.Content.Find.Execute("6+4=", $true, $false, $false, $false, $false, $true, 1, $false, "37+16=", 2) | Out-Null
$d.Content.Find.Execute("72-38=", $true, $false, $false, $false, $false, $true, 1, $false, "62-52=", 2) | Out-Null

# --- Step 2: insert 4 new rows before (original) row 13, and fill their text ---
# NOTE: Rows.Add(beforeRow) always inserts immediately above beforeRow, so
# repeated Add calls build the block in reverse order. We add the rows in
# reverse order here so the final on-page order matches the target.
$refRow = $t.Rows.Item(13)
$newRowsContent = @(
    @("50-34=", "43+8=", "90-45=", "75-1=", "44+3="),
    @("29+23=", "51+28=", "71-57=", "98+1=", "18+2="),
    @("82-61=", "35-1=", "10+4=", "14+54=", "38-14="),
    @("65+22=", "14+10=", "18-15=", "98-5=", "72-71=")
)
for ($r = $newRowsContent.Count - 1; $r -ge 0; $r--) {
    $rowVals = $newRowsContent[$r]
    $newRow = $t.Rows.Add($refRow)
    for ($i = 1; $i -le $rowVals.Count; $i++) {
        $newRow.Cells.Item($i).Range.Text = $rowVals[$i - 1]
    }
}

# --- Step 3: delete the last 4 rows (original rows 17-20) ---
for ($k = 0; $k -lt 4; $k++) {
    $lastIndex = $t.Rows.Count
    $t.Rows.Item($lastIndex).Delete()
}

Write-Output "Rows after edit: $($t.Rows.Count)"
